# Tenacious Estimation - WBS: update per-task hour estimates,
# add ESTIMATED TIME / BUFFER HOURS / TOTAL summary rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated "Estimation in Hrs" (column D) values -------------------------
$ws.Range("D2").Value = 1
$ws.Range("D4").Value = 0.5
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 0.5
$ws.Range("D8").Value = 1.5
$ws.Range("D10").Value = 1
$ws.Range("D12").Value = 0.5
$ws.Range("D13").Value = 1
$ws.Range("D14").Value = 0.5
$ws.Range("D15").Value = 0.5
$ws.Range("D16").Value = 1
$ws.Range("D17").Value = 0.5
$ws.Range("D18").Value = 1
$ws.Range("D48").Value = 0.5
$ws.Range("D49").Value = 1
$ws.Range("D51").Value = 1
$ws.Range("D52").Value = 1
$ws.Range("D54").Value = 1.5
$ws.Range("D58").Value = 2
$ws.Range("D60").Value = 2
$ws.Range("D82").Value = 1
$ws.Range("D83").Value = 1
$ws.Range("D84").Value = 0.2
$ws.Range("D85").Value = 1
$ws.Range("D86").Value = 1
$ws.Range("D88").Value = 1
$ws.Range("D95").Value = 1
$ws.Range("D97").Value = 0.5
$ws.Range("D98").Value = 2
$ws.Range("D99").Value = 2
$ws.Range("D104").Value = 1
$ws.Range("D105").Value = 1
$ws.Range("D106").Value = 1.5
$ws.Range("D108").Value = 0.5
$ws.Range("D109").Value = 0.5
$ws.Range("D114").Value = 0.5
$ws.Range("D118").Value = 1
$ws.Range("D119").Value = 1
$ws.Range("D121").Value = 1
$ws.Range("D122").Value = 1
$ws.Range("D124").Value = 1
$ws.Range("D127").Value = 0.5
$ws.Range("D128").Value = 1
$ws.Range("D129").Value = 1
$ws.Range("D130").Value = 1
$ws.Range("D131").Value = 1
$ws.Range("D132").Value = 1
$ws.Range("D134").Value = 1
$ws.Range("D135").Value = 0.5
$ws.Range("D136").Value = 0.5
$ws.Range("D137").Value = 0.5
$ws.Range("D141").Value = 1
$ws.Range("D142").Value = 1
$ws.Range("D143").Value = 0.2
$ws.Range("D144").Value = 1
$ws.Range("D145").Value = 1
$ws.Range("D147").Value = 1
$ws.Range("D149").Value = 0.5
$ws.Range("D151").Value = 1
$ws.Range("D153").Value = 0.5
$ws.Range("D154").Value = 2
$ws.Range("D155").Value = 2
$ws.Range("D158").Value = 1
$ws.Range("D160").Value = 1
$ws.Range("D164").Value = 1
$ws.Range("D165").Value = 1.5
$ws.Range("D167").Value = 0.5
$ws.Range("D168").Value = 0.5
$ws.Range("D173").Value = 0.5
$ws.Range("D174").Value = 1
$ws.Range("D175").Value = 1
$ws.Range("D176").Value = 1
$ws.Range("D177").Value = 1
$ws.Range("D178").Value = 1
$ws.Range("D180").Value = 1
$ws.Range("D181").Value = 1
$ws.Range("D183").Value = 0.5
$ws.Range("D221").Value = 1
$ws.Range("D222").Value = 1
$ws.Range("D224").Value = 1
$ws.Range("D225").Value = 1
$ws.Range("D227").Value = 1
$ws.Range("D228").Value = 1
$ws.Range("D230").Value = 0.5
$ws.Range("D336").Value = 1
$ws.Range("D338").Value = 0.5
$ws.Range("D339").Value = 1.5
$ws.Range("D340").Value = 1
$ws.Range("D342").Value = 0.5
$ws.Range("D345").Value = 1
$ws.Range("D346").Value = 1
$ws.Range("D347").Value = 1
$ws.Range("D349").Value = 1.5
$ws.Range("D350").Value = 1.5
$ws.Range("D352").Value = 1
$ws.Range("D355").Value = 1
$ws.Range("D356").Value = 1
$ws.Range("D357").Value = 1.5
$ws.Range("D358").Value = 0
$ws.Range("D359").Value = 0.5
$ws.Range("D360").Value = 0.5
$ws.Range("D365").Value = 0.5
$ws.Range("D369").Value = 1
$ws.Range("D371").Value = 0.5
$ws.Range("D372").Value = 1.5
$ws.Range("D373").Value = 1
$ws.Range("D374").Value = 1
$ws.Range("D375").Value = 0.5
$ws.Range("D378").Value = 1
$ws.Range("D379").Value = 1
$ws.Range("D380").Value = 1.5
$ws.Range("D381").Value = 0
$ws.Range("D382").Value = 0.5
$ws.Range("D383").Value = 0.5
$ws.Range("D387").Value = 0.5
$ws.Range("D388").Value = 0.5
$ws.Range("D389").Value = 1.5
$ws.Range("D392").Value = 2
$ws.Range("D393").Value = 2
$ws.Range("D396").Value = 1.5
$ws.Range("D399").Value = 1.5
$ws.Range("D402").Value = 2
$ws.Range("D403").Value = 2
$ws.Range("D406").Value = 1.5
$ws.Range("D412").Value = 1
$ws.Range("D415").Value = 2
$ws.Range("D416").Value = 2
$ws.Range("D418").Value = 1
$ws.Range("D421").Value = 1.5
$ws.Range("D422").Value = 1
$ws.Range("D425").Value = 2
$ws.Range("D426").Value = 2
$ws.Range("D428").Value = 1
$ws.Range("D431").Value = 1.5
$ws.Range("D435").Value = 2
$ws.Range("D436").Value = 2
$ws.Range("D441").Value = 1.5
$ws.Range("D11").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("D50").ClearContents()
$ws.Range("D107").ClearContents()
$ws.Range("D166").ClearContents()

# --- New summary rows below the grand total --------------------------------
# Row 444 already holds =SUM(D2:D441) in D444; add its label.
$ws.Range("C444").Value = "ESTIMATED TIME"

# Row 445: buffer hours.
$ws.Range("C445").Value = "BUFFER HOURS "
$ws.Range("D445").Value = 20

# Row 448: grand total (estimated time + buffer hours), centered + bold-free
# style matching the new centered cellXf.
$ws.Range("C448").Value = "TOTAL"
$ws.Range("D448").Value = 606.4
$ws.Range("C448:D448").HorizontalAlignment = -4108

# --- Restore the view state (scroll position / active selection) ----------
[void]$ws.Range("A424").Select()
$excel.ActiveWindow.ScrollRow = 424
[void]$ws.Range("C448").Select()
